$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(newC, newD)  (row 82 only has C; D is $null meaning "skip")
$updates = @{
    2 = @(96, 71)
    3 = @(100, 79)
    4 = @(55, 38)
    5 = @(65, 44.5)
    6 = @(5, 76)
    7 = @(78, 54.5)
    8 = @(104, 84)
    9 = @(32, 67)
    10 = @(27, 58.5)
    11 = @(163, 161.5)
    12 = @(62, 40.5)
    13 = @(149, 129.5)
    14 = @(113, 93.5)
    15 = @(57, 33.5)
    16 = @(185, 185)
    17 = @(15, 75.5)
    18 = @(16, 74)
    19 = @(14, 75.5)
    20 = @(53, 30)
    21 = @(26, 26)
    22 = @(95, 80.5)
    23 = @(145, 129)
    24 = @(61, 39.5)
    25 = @(17, 69)
    26 = @(38, 63.5)
    27 = @(117, 93.5)
    28 = @(28, 19.5)
    29 = @(4, 77)
    30 = @(74, 50.5)
    31 = @(105, 99)
    33 = @(107, 81.5)
    34 = @(93, 71)
    35 = @(205, 207)
    36 = @(272, 267.5)
    37 = @(11, 73)
    38 = @(152, 130.5)
    39 = @(136, 124)
    40 = @(108, 85)
    41 = @(66, 46)
    42 = @(18, 79)
    43 = @(80, 58)
    44 = @(31, 64)
    45 = @(48, 25.5)
    46 = @(247, 241.5)
    47 = @(9, 74)
    48 = @(51, 33.5)
    49 = @(203, 202.5)
    50 = @(218, 218.5)
    51 = @(46, 25.5)
    52 = @(137, 125.5)
    53 = @(6, 76)
    54 = @(8, 80)
    55 = @(91, 67.5)
    56 = @(56, 44)
    57 = @(7, 78)
    58 = @(88, 70.5)
    59 = @(45, 33.5)
    60 = @(69, 47)
    61 = @(85, 69.5)
    62 = @(72, 55)
    63 = @(289, 282)
    64 = @(60, 37.5)
    65 = @(106, 83)
    66 = @(79, 43.5)
    67 = @(82, 61.5)
    68 = @(146, 153.5)
    69 = @(47, 35)
    70 = @(118, 94.5)
    71 = @(77, 60)
    72 = @(43, 69.5)
    73 = @(50, 31)
    74 = @(180, 179)
    75 = @(102, 98.5)
    76 = @(87, 69.5)
    77 = @(75, 52)
    78 = @(101, 79)
    79 = @(193, 184)
    80 = @(120, 102.5)
    81 = @(54, 33.5)
    82 = @(87.3, $null)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals[0]) { $ws.Cells.Item($row, 3).Value = $vals[0] }
    if ($null -ne $vals[1]) { $ws.Cells.Item($row, 4).Value = $vals[1] }
}
